# Increment the "想去人数" (F column) counts by 1 for specific rows
# in both the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")
$rows = @(2, 3, 7, 8)

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($r in $rows) {
        $cell = $ws.Cells.Item($r, 6)  # Column F = 6
        $cell.Value2 = $cell.Value2 + 1
    }
}
